$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'51.603.03"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.77%  "
$ws.Cells.Item(3, 4).Value = "'2.793.71"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +1.40%  "
$ws.Cells.Item(4, 5).Value = "  -0.06%  "
$ws.Cells.Item(5, 4).Value = "'352.43"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.80%  "
$ws.Cells.Item(6, 4).Value = "'111.20"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +3.18%  "
$ws.Cells.Item(7, 4).Value = "'0.554"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.50%  "
$ws.Cells.Item(8, 5).Value = "  -0.06%  "
$ws.Cells.Item(9, 4).Value = "'0.627"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +8.27%  "
$ws.Cells.Item(10, 4).Value = "'39.99"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +1.96%  "
$ws.Cells.Item(11, 5).Value = "  -1.15%  "
$ws.Cells.Item(12, 4).Value = "'0.0834"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.11%  "
$ws.Cells.Item(13, 4).Value = "'19.89"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +2.03%  "
$ws.Cells.Item(14, 4).Value = "'7.73"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +3.08%  "
$ws.Cells.Item(15, 4).Value = "'3.233.04"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.04%  "
$ws.Cells.Item(16, 4).Value = "'2.799.84"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +1.41%  "
$ws.Cells.Item(17, 4).Value = "'0.941"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +2.28%  "
$ws.Cells.Item(18, 4).Value = "'51.598.14"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.51%  "
$ws.Cells.Item(19, 4).Value = "'7.56"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.58%  "
$ws.Cells.Item(20, 4).Value = "'3.19"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +6.57%  "
$ws.Cells.Item(21, 4).Value = "'13.49"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +3.69%  "
$ws.Cells.Item(22, 4).Value = "'0.0₃0967"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.92%  "
$ws.Cells.Item(23, 4).Value = "'70.15"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.94%  "
$ws.Cells.Item(24, 4).Value = "'267.04"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.43%  "
$ws.Cells.Item(25, 5).Value = "  +0.45%  "
$ws.Cells.Item(26, 5).Value = "  +0.08%  "
$ws.Cells.Item(27, 4).Value = "'26.02"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.21%  "
$ws.Cells.Item(28, 4).Value = "'0.160"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.35%  "
$ws.Cells.Item(29, 4).Value = "'38.95"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +11.74%  "
$ws.Cells.Item(30, 4).Value = "'10.31"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.56%  "
$ws.Cells.Item(31, 5).Value = "  +0.47%  "
$ws.Cells.Item(32, 4).Value = "'52.61"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +2.37%  "
$ws.Cells.Item(33, 4).Value = "'6.09"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.96%  "
$ws.Cells.Item(34, 4).Value = "'0.0452"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +3.31%  "
$ws.Cells.Item(35, 4).Value = "'0.0889"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +6.88%  "
$ws.Cells.Item(36, 4).Value = "'5.56"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +8.64%  "
$ws.Cells.Item(37, 5).Value = "  -0.22%  "
$ws.Cells.Item(38, 4).Value = "'18.75"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.88%  "
$ws.Cells.Item(39, 5).Value = "  +3.50%  "
$ws.Cells.Item(40, 4).Value = "'3.15"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +0.72%  "
$ws.Cells.Item(41, 5).Value = "  +1.43%  "
$ws.Cells.Item(42, 5).Value = "  +0.67%  "
$ws.Cells.Item(43, 5).Value = "  +1.28%  "
$ws.Cells.Item(44, 4).Value = "'121.17"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.73%  "
$ws.Cells.Item(45, 4).Value = "'21.79"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.29%  "
$ws.Cells.Item(46, 4).Value = "'2.44"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +5.91%  "
$ws.Cells.Item(47, 5).Value = "  +5.48%  "
$ws.Cells.Item(48, 4).Value = "'2.105.35"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.37%  "
$ws.Cells.Item(49, 4).Value = "'0.960"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +4.63%  "
$ws.Cells.Item(50, 4).Value = "'5.46"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.68%  "
$ws.Cells.Item(51, 5).Value = "  +6.55%  "
